# Apply updated experiment values ("add almost all lineal")
$wb = $excel.ActiveWorkbook

# --- Restricciones_del_follower ---
$ws = $wb.Worksheets.Item("Restricciones_del_follower")

$ws.Range("A2").Value = "0.01514095315217201y_1 + 0.9127063484843593y_2"
$ws.Range("B2").Value = "'2.2456891574731968"
$ws.Range("D2").Value = "'0.0866877650392671"
$ws.Range("E2").Value = "'0.16392888961395025"
$ws.Range("F2").Value = "'0.5256365753464083"

$ws.Range("A3").Value = "-4 + 0.4641757826909333y_1 - 0.48175592097931874y_2"
$ws.Range("B3").Value = "'-3.115964380723809"
$ws.Range("D3").Value = "'0.9648587319705634"
$ws.Range("E3").Value = "'0.9400307632963608"
$ws.Range("F3").Value = "'0.019517476239305102"

$ws.Range("A4").Value = "-16 - 2x - 1.1497178456235067y_1 + 2.06720325974616y_2"
$ws.Range("B4").Value = "'-27.854581994052488"
$ws.Range("D4").Value = "'0.9761226555169311"
$ws.Range("E4").Value = "'0.6645676367532283"
$ws.Range("F4").Value = "'0.4658746234119081"

$ws.Range("A5").Value = "-48 + 8x + 1.1821095518371985y_1 + 0.1637334633829287y_2"
$ws.Range("B5").Value = "'4.57854615527611"
$ws.Range("D5").Value = "'0.8143958706897286"
$ws.Range("E5").Value = "'0.21445557599747056"
$ws.Range("F5").Value = "'0.8360641848199136"

$ws.Range("A6").Value = "12 - 2x - 0.3895597410918963y_1 + 1.4479359182545817y_2"
$ws.Range("B6").Value = "'1.9983186674708953"
$ws.Range("D6").Value = "'0.23927405565041526"
$ws.Range("E6").Value = "'0.9861708887310502"
$ws.Range("F6").Value = "'0.5737634833120397"

# --- Punto_modificado ---
$ws = $wb.Worksheets.Item("Punto_modificado")
$ws.Range("A2").Value = "'5.875840352759835"
$ws.Range("B2").Value = "'4.382729079133727"
$ws.Range("C2").Value = "'2.387767396848251"

# --- Vector_bf ---
# NOTE: worksheet name lookup by "Item(name)" is case-insensitive, and this
# workbook has two sheets whose names differ only by case ("Vector_bf" vs
# "Vector_BF"). Use the stable 1-based sheet index instead so each one is
# addressed unambiguously.
$ws = $wb.Worksheets.Item(5)   # Vector_bf
$ws.Range("A2").Value = "'0.8035954455273624"
$ws.Range("A3").Value = "'-2.111935357803988"

# --- Vector_BF ---
$ws = $wb.Worksheets.Item(6)   # Vector_BF
$ws.Range("A2").Value = "'2.5858324429887922"
$ws.Range("A3").Value = "'3.4559062079095257"
$ws.Range("A4").Value = "'-4.533575742678566"

# --- Vector_Alpha (plain numeric cells) ---
$ws = $wb.Worksheets.Item("Vector_Alpha")
$ws.Range("A2").Value = 0.6289406846527845
$ws.Range("A3").Value = 0.5654763054531662
